# Update the "generator" linear example (ex9.1.2, alpha_zero) with the
# values produced by a fresh run of the random generator.
#
# All of the numeric-looking values in this workbook were originally
# authored as TEXT (shared strings with t="s"), not as numbers, so we
# force the target ranges to Text format before writing the new values;
# otherwise Excel would silently re-interpret "0.824417605514952" etc.
# as a real number.

$wb = $excel.ActiveWorkbook

# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively here, so we address every sheet by its (unambiguous)
# 1-based tab position instead:
#   1 Funciones_Objetivo        5 Vector_bf
#   2 Restricciones_del_lider   6 Vector_BF
#   3 Restricciones_del_follower 7 Vector_Alpha
#   4 Punto_modificado

# ---------------------------------------------------------------------
# Sheet: Restricciones_del_follower
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

# NOTE: multi-area ranges ("B2:B5,D2:F5") only apply NumberFormat to the
# first area in this runtime, so the format is applied one contiguous
# block at a time instead.
$ws.Range("B2:B5").NumberFormat = "@"
$ws.Range("D2:F5").NumberFormat = "@"

$ws.Range("B2").Value = "-4.374623078112156"
$ws.Range("D2").Value = "0.33468162538227564"
$ws.Range("E2").Value = "0.5786326093455703"
$ws.Range("F2").Value = "0.824417605514952"

$ws.Range("B3").Value = "-2.8102693382873367"
$ws.Range("D3").Value = "0.9092567913461869"
$ws.Range("E3").Value = "0.6604724023704651"
$ws.Range("F3").Value = "0.399906499902034"

$ws.Range("B4").Value = "0.9341385726238034"
$ws.Range("D4").Value = "0.7906785535517057"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0.389168975639538"

$ws.Range("B5").Value = "0.36494658748581443"
$ws.Range("D5").Value = "0.5618257705012442"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0.7234541633898545"

# ---------------------------------------------------------------------
# Sheet: Punto_modificado
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)

$rngPunto = $ws.Range("A2:B2")
$rngPunto.NumberFormat = "@"

$ws.Range("A2").Value = "4.184892416399492"
$ws.Range("B2").Value = "4.374623078112156"

# ---------------------------------------------------------------------
# Sheet: Vector_bf
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

$rngBf = $ws.Range("A2")
$rngBf.NumberFormat = "@"

$ws.Range("A2").Value = "-2.5941065025660786"

# ---------------------------------------------------------------------
# Sheet: Vector_BF
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

$rngBF = $ws.Range("A2:A3")
$rngBF.NumberFormat = "@"

$ws.Range("A2").Value = "1.6604724023704651"
$ws.Range("A3").Value = "2.918160206975105"
